$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Fitness) rows 2-96 currently hold 7310; update them to 7293
# to match the rest of the column (rows 97-252 already hold 7293).
$ws.Range("C2:C96").Value = 7293
